$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (date format hint changed from MM/dd/yyyy to yyyy/MM/dd) ---
$ws.Range("B1").Value = "Trip Pick up Date Start yyyy/MM/dd (Text Format)   تاريخ بداية الرحلة"
$ws.Range("C1").Value = "Trip Pick up Date End yyyy/MM/dd (Text Format)   تاريخ نهاية الرحلة"

# --- Columns B and C should be plain Text format (no more date-picker format) ---
$ws.Columns("B").NumberFormat = "@"
$ws.Columns("C").NumberFormat = "@"

# --- C3/C4 previously carried an explicit date number format (m/d/yyyy); clear it
#     back to plain text so it matches the rest of the column ---
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"

# --- Selection: row 2 (the first data-entry row) is now selected in its entirety,
#     and the view no longer has a frozen/scrolled topLeftCell ---
$ws.Rows(2).Select()

Write-Host "done"
